$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The contacts sheet had two rows pointing at the same phone number / group
# ("Jasper Barcelona" and "Leanza Etorma", both in "New Group"). The fix
# collapses this back down to a single, corrected contact row: keep the
# phone number in A1, but put "Leanza Etorma" in B1 and re-assign her to
# "Sample Group" in C1 - then drop the now-redundant second row entirely.
$ws.Range("B1").Value = "Leanza Etorma"
$ws.Range("C1").Value = "Sample Group"
$ws.Rows.Item(2).Delete()

# Columns re-flow (narrower name column, a dedicated width for the group
# column) once the sheet is down to a single row of data.
$ws.Columns.Item(1).ColumnWidth = 16.666666666666668
$ws.Columns.Item(2).ColumnWidth = 28.666666666666668
$ws.Columns.Item(3).ColumnWidth = 11.833333333333334

# Selection moves to the group's name cell.
$ws.Range("B1").Select()
